# "price, seo and adsense update" - refresh the Power Query-backed "Table 0" price
# table (Gasoline_Price / Diesel_Price / LPG_Price columns) with the latest fetched
# values for every country row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Table 0" - the query result sheet (tabSelected="1")

# Row 2: Austria
$ws.Range("B2").Value = "196,96 "
$ws.Range("C2").Value = "198,72 "
$ws.Range("D2").Value = "138,85 "
# Row 3: Bosnia and Herzegovina
$ws.Range("B3").Value = "147,21 "
$ws.Range("C3").Value = "150,80 "
$ws.Range("D3").Value = "81,98 "
# Row 4: Bulgaria
$ws.Range("B4").Value = "157,53 "
$ws.Range("C4").Value = "158,73 "
$ws.Range("D4").Value = "77,87 "
# Row 5: Switzerland
$ws.Range("B5").Value = "210,43 "
$ws.Range("C5").Value = "225,14 "
$ws.Range("D5").Value = "118,06 "
# Row 6: Czech Republic
$ws.Range("B6").Value = "179,47 "
$ws.Range("C6").Value = "172,06 "
$ws.Range("D6").Value = "83,59 "
# Row 7: Germany
$ws.Range("B7").Value = "212,19 "
$ws.Range("C7").Value = "197,78 "
$ws.Range("D7").Value = "128,65 "
# Row 8: Estonia
$ws.Range("B8").Value = "192,39 "
$ws.Range("C8").Value = "181,85 "
$ws.Range("D8").Value = "121,97 "
# Row 9: Spain
$ws.Range("B9").Value = "190,17 "
$ws.Range("C9").Value = "183,49 "
$ws.Range("D9").Value = "114,94 "
# Row 10: France
$ws.Range("B10").Value = "214,66 "
$ws.Range("C10").Value = "206,81 "
$ws.Range("D10").Value = "121,74 "
# Row 11: United Kingdom
$ws.Range("B11").Value = "216,93 "
$ws.Range("C11").Value = "239,26 "
$ws.Range("D11").Value = "133,13 "
# Row 12: Greece
$ws.Range("B12").Value = "215,59 "
$ws.Range("C12").Value = "187,47 "
$ws.Range("D12").Value = "115,41 "
# Row 13: Croatia
$ws.Range("B13").Value = "178,10 "
$ws.Range("C13").Value = "168,72 "
$ws.Range("D13").Value = "98,42 "
# Row 14: Hungary
$ws.Range("B14").Value = "183,33 "
$ws.Range("C14").Value = "185,23 "
$ws.Range("D14").Value = "106,73 "
# Row 15: Ireland
$ws.Range("B15").Value = "211,26 "
$ws.Range("C15").Value = "209,27 "
$ws.Range("D15").Value = "117,17 "
# Row 16: Italy
$ws.Range("B16").Value = "224,15 "
$ws.Range("C16").Value = "213,13 "
$ws.Range("D16").Value = "86,24 "
# Row 17: Lithuania
$ws.Range("B17").Value = "177,04 "
$ws.Range("C17").Value = "189,82 "
$ws.Range("D17").Value = "93,97 "
# Row 18: Latvia
$ws.Range("B18").Value = "191,46 "
$ws.Range("C18").Value = "189,00 "
$ws.Range("D18").Value = "109,79 "
# Row 19: North Macedonia
$ws.Range("B19").Value = "150,39 "
$ws.Range("C19").Value = "139,92 "
$ws.Range("D19").Value = "87,57 "
# Row 20: Netherlands
$ws.Range("B20").Value = "246,53 "
$ws.Range("C20").Value = "220,51 "
$ws.Range("D20").Value = "106,98 "
# Row 21: Poland
$ws.Range("B21").Value = "177,20 "
$ws.Range("C21").Value = "180,03 "
$ws.Range("D21").Value = "94,68 "
# Row 22: Portugal
$ws.Range("B22").Value = "214,89 "
$ws.Range("C22").Value = "208,09 "
$ws.Range("D22").Value = "113,07 "
# Row 23: Romania
$ws.Range("B23").Value = "178,20 "
$ws.Range("C23").Value = "183,85 "
$ws.Range("D23").Value = "84,28 "
# Row 24: Serbia
$ws.Range("B24").Value = "187,98 "
$ws.Range("C24").Value = "203,01 "
$ws.Range("D24").Value = "109,04 "
# Row 25: Slovenia
$ws.Range("B25").Value = "183,02 "
$ws.Range("C25").Value = "189,58 "
$ws.Range("D25").Value = "107,68 "
# Row 26: Slovakia
$ws.Range("B26").Value = "190,52 "
$ws.Range("C26").Value = "186,89 "
$ws.Range("D26").Value = "91,16 "
# Row 27: Turkey
$ws.Range("B27").Value = "145,53 "
$ws.Range("C27").Value = "150,72 "
$ws.Range("D27").Value = "83,86 "
